# Add a new "Team" column (AF) to the Mikel Arteta sheet.
# Column AF records the coach's club ("Arsenal") for every match row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell AF1 - set the text first, then copy the header formatting
# (bold font, thin border, centered/top-aligned) from the neighbouring
# ProbA header (AE1) so the new header matches the existing ones.
$ws.Range("AF1").Value = "Team"
$ws.Range("AE1").Copy()
$ws.Range("AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-59: every match in this sheet was coached for Arsenal.
$lastRow = 59
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 32).Value = "Arsenal"
}
